# Update countries & provincias Spain
# Applies the data refresh for "Datos actualizados a 21 de Mayo de 2020":
#  - Updated case counts for several existing countries (rows keep their
#    country but the day's figures changed).
#  - Re-sorted entries for "Guinea Ecuatorial" (moved above Paraguay /
#    Zambia) and "Mauritania" (moved above Brunei / Mongolia) because
#    their case totals overtook the countries that used to be ranked
#    above them; the rows in between simply shift down one position.
#  - Refreshed the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 00:35"

# --- Estados Unidos (row 4) --------------------------------------------
$ws.Range("B4").Value = 1589834
$ws.Range("C4").Value = 19251
$ws.Range("D4").Value = 369362
$ws.Range("E4").Value = 1125642
$ws.Range("G4").Value = 1297
$ws.Range("H4").Value = 94830

# --- Brasil (row 6) -----------------------------------------------------
$ws.Range("B6").Value = 291579
$ws.Range("C6").Value = 19694
$ws.Range("D6").Value = 116683
$ws.Range("E6").Value = 156037
$ws.Range("G6").Value = 876
$ws.Range("H6").Value = 18859

# --- row 11 ---------------------------------------------------------
$ws.Range("B11").Value = 178531
$ws.Range("C11").Value = 704
$ws.Range("E11").Value = 13361
$ws.Range("G11").Value = 77
$ws.Range("H11").Value = 8270

# --- Argentina (row 51) -------------------------------------------------
$ws.Range("B51").Value = 9283
$ws.Range("C51").Value = 474
$ws.Range("E51").Value = 5947
$ws.Range("G51").Value = 10
$ws.Range("H51").Value = 403

# --- Chequia (row 52) ----------------------------------------------------
$ws.Range("B52").Value = 8721
$ws.Range("C52").Value = 74
$ws.Range("D52").Value = 5830
$ws.Range("E52").Value = 2587

# --- Guinea Ecuatorial / Paraguay / Zambia re-sort (rows 114-116) -------
$ws.Range("A114").Value = "Guinea Ecuatorial"
$ws.Range("B114").Value = 890
$ws.Range("C114").Value = 65
$ws.Range("D114").Value = 22
$ws.Range("E114").Value = 861
$ws.Range("H114").Value = 7

$ws.Range("A115").Value = "Paraguay"
$ws.Range("B115").Value = 833
$ws.Range("C115").Value = 4
$ws.Range("D115").Value = 242
$ws.Range("E115").Value = 580
$ws.Range("H115").Value = 11

$ws.Range("A116").Value = "Zambia"
$ws.Range("B116").Value = 832
$ws.Range("C116").Value = 60
$ws.Range("D116").Value = 197
$ws.Range("E116").Value = 628

# --- Mauritania / Brunei / Mongolia re-sort (rows 160-162) --------------
$ws.Range("A160").Value = "Mauritania"
$ws.Range("C160").Value = 10
$ws.Range("D160").Value = 7
$ws.Range("E160").Value = 130
$ws.Range("H160").Value = 4

$ws.Range("A161").Value = "Brunei"
$ws.Range("B161").Value = 141
$ws.Range("D161").Value = 136
$ws.Range("E161").Value = 4
$ws.Range("H161").Value = 1

$ws.Range("A162").Value = "Mongolia"
$ws.Range("B162").Value = 140
$ws.Range("D162").Value = 26
$ws.Range("E162").Value = 114

# --- row 164 --------------------------------------------------------
$ws.Range("D164").Value = 47
$ws.Range("E164").Value = 68
